$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings are preserved as text (matching the
# source data which stores every Price/Volume cell as a literal string).
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

$ws.Range("D2").Value = "61.915.70"
$ws.Range("E2").Value = "  -1.16%  "
$ws.Range("D3").Value = "3.407.84"
$ws.Range("E3").Value = "  -0.62%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("E5").Value = "  +0.41%  "
$ws.Range("D6").Value = "128.67"
$ws.Range("E6").Value = "  -1.41%  "
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "0.733"
$ws.Range("E9").Value = "  +5.66%  "
$ws.Range("E10").Value = "  +1.98%  "
$ws.Range("D11").Value = "42.93"
$ws.Range("E11").Value = "  +1.85%  "
$ws.Range("E12").Value = "  +36.12%  "
$ws.Range("D13").Value = "9.31"
$ws.Range("E13").Value = "  +9.89%  "
$ws.Range("E14").Value = "  -0.38%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "21.30"
$ws.Range("E15").Value = "  +7.46%  "
$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").Value = "3.949.32"
$ws.Range("E16").Value = "  -0.75%  "
$ws.Range("D17").Value = "3.412.38"
$ws.Range("E17").Value = "  -0.56%  "
$ws.Range("D18").Value = "12.56"
$ws.Range("E18").Value = "  +8.86%  "
$ws.Range("E19").Value = "  +6.82%  "
$ws.Range("D20").Value = "61.937.60"
$ws.Range("E20").Value = "  -1.10%  "
$ws.Range("D21").Value = "447.39"
$ws.Range("E21").Value = "  +42.04%  "
$ws.Range("D22").Value = "92.23"
$ws.Range("E22").Value = "  +9.00%  "
$ws.Range("E23").Value = "  +0.51%  "
$ws.Range("D24").Value = "13.17"
$ws.Range("E24").Value = "  +2.61%  "
$ws.Range("E25").Value = "  +3.42%  "
$ws.Range("D26").Value = "9.38"
$ws.Range("E26").Value = "  +14.76%  "
$ws.Range("D27").Value = "33.15"
$ws.Range("E27").Value = "  +11.30%  "
$ws.Range("E28").Value = "  +0.71%  "
$ws.Range("D29").Value = "7.76"
$ws.Range("E29").Value = "  -0.75%  "
$ws.Range("E30").Value = "  -2.53%  "
$ws.Range("D31").Value = "11.97"
$ws.Range("E31").Value = "  +4.94%  "
$ws.Range("E32").Value = "  -1.42%  "
$ws.Range("E33").Value = "  -0.31%  "
$ws.Range("D34").Value = "42.73"
$ws.Range("E34").Value = "  -3.67%  "
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("E36").Value = "  +4.09%  "
$ws.Range("D37").Value = "53.82"
$ws.Range("E37").Value = "  +3.78%  "
$ws.Range("E39").Value = "  +1.63%  "
$ws.Range("E40").Value = "  +7.59%  "
$ws.Range("D41").Value = "2.96"
$ws.Range("E41").Value = "  -0.23%  "
$ws.Range("D42").Value = "0.320"
$ws.Range("E42").Value = "  -0.71%  "
$ws.Range("D43").Value = "143.61"
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("D44").Value = "4.31"
$ws.Range("E44").Value = "  +9.89%  "
$ws.Range("E45").Value = "  +15.11%  "
$ws.Range("E46").Value = "  +0.84%  "
$ws.Range("D47").Value = "16.61"
$ws.Range("E47").Value = "  -1.77%  "
$ws.Range("D48").Value = "0.148"
$ws.Range("E48").Value = "  +22.28%  "
$ws.Range("D49").Value = "22.34"
$ws.Range("E49").Value = "  +4.29%  "
$ws.Range("E50").Value = "  +7.58%  "
$ws.Range("D51").Value = "3.749.82"
$ws.Range("E51").Value = "  -0.71%  "
